# Insert a new weekly record for "Espinaca" (Terminal La Palmera de La Serena)
# at row 295, shifting the existing rows 295:318 down to 296:319.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 295:318 down one row, leaving a fresh blank row 295.
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(295, 1).Value = 8
$ws.Cells.Item(295, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(295, 3).Value = "Coquimbo"
$ws.Cells.Item(295, 4).Value = 44826
$ws.Cells.Item(295, 5).Value = 4
$ws.Cells.Item(295, 6).Value = 100112012
$ws.Cells.Item(295, 7).Value = "Espinaca"
$ws.Cells.Item(295, 8).Value = "Sin especificar"
$ws.Cells.Item(295, 9).Value = "Primera"
$ws.Cells.Item(295, 10).Value = 2800
$ws.Cells.Item(295, 11).Value = 450
$ws.Cells.Item(295, 12).Value = 500
$ws.Cells.Item(295, 13).Value = 475
$ws.Cells.Item(295, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(295, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(295, 16).Value = 950
$ws.Cells.Item(295, 17).Value = 0.5
$ws.Cells.Item(295, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of column D.
$ws.Cells.Item(295, 4).NumberFormat = $ws.Cells.Item(296, 4).NumberFormat
